# Scheduled runner update: refresh profit-calc columns (H:N) across the
# per-job Leve sheets with newly pulled market-board averages.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1704.174
$ws.Range("I70").Value = 1498.75
$ws.Range("K70").Value = 4496.25
$ws.Range("M70").Value = -4226.25

$ws.Range("H73").Value = 1704.174
$ws.Range("I73").Value = 1498.75
$ws.Range("K73").Value = 4496.25
$ws.Range("M73").Value = -3560.25

$ws.Range("H129").Value = 1376.5769
$ws.Range("J129").Value = 869.3043
$ws.Range("L129").Value = 2607.9129
$ws.Range("N129").Value = -12607.9129

$ws.Range("H132").Value = 33176.805
$ws.Range("I132").Value = 5122.619
$ws.Range("J132").Value = 92090.60000000001
$ws.Range("K132").Value = 15367.857
$ws.Range("L132").Value = 276271.8
$ws.Range("M132").Value = -12837.857
$ws.Range("N132").Value = -281331.8

$ws.Range("H137").Value = 12065.048
$ws.Range("I137").Value = 16470.5
$ws.Range("J137").Value = 9354
$ws.Range("K137").Value = 49411.5
$ws.Range("L137").Value = 28062
$ws.Range("M137").Value = -46861.5
$ws.Range("N137").Value = -33162

$ws.Range("H138").Value = 2223.85
$ws.Range("I138").Value = 1153.8
$ws.Range("J138").Value = 2988.1714
$ws.Range("K138").Value = 3461.4
$ws.Range("L138").Value = 8964.514200000001
$ws.Range("M138").Value = 1678.6
$ws.Range("N138").Value = -19244.5142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3775
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 6100
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 6100
$ws.Range("M2").Value = -2887
$ws.Range("N2").Value = -6326

$ws.Range("H32").Value = 15449.667
$ws.Range("I32").Value = 16160.243
$ws.Range("K32").Value = 16160.243
$ws.Range("M32").Value = -15873.243

$ws.Range("H45").Value = 800
$ws.Range("I45").Value = 800
$ws.Range("K45").Value = 800
$ws.Range("M45").Value = -423

$ws.Range("H74").Value = 1606.2245
$ws.Range("I74").Value = 1392.9688
$ws.Range("J74").Value = 2007.6471
$ws.Range("K74").Value = 1392.9688
$ws.Range("L74").Value = 2007.6471
$ws.Range("M74").Value = -518.9688000000001
$ws.Range("N74").Value = -3755.6471

$ws.Range("H77").Value = 1606.2245
$ws.Range("I77").Value = 1392.9688
$ws.Range("J77").Value = 2007.6471
$ws.Range("K77").Value = 6964.844000000001
$ws.Range("L77").Value = 10038.2355
$ws.Range("M77").Value = -2596.844000000001
$ws.Range("N77").Value = -18774.2355

$ws.Range("H110").Value = 5180.7334
$ws.Range("I110").Value = 2190.111
$ws.Range("J110").Value = 9666.666999999999
$ws.Range("K110").Value = 2190.111
$ws.Range("L110").Value = 9666.666999999999
$ws.Range("M110").Value = -145.1109999999999
$ws.Range("N110").Value = -13756.667

$ws.Range("H116").Value = 3775
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 6100
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 6100
$ws.Range("M116").Value = -706
$ws.Range("N116").Value = -10688

$ws.Range("H132").Value = 15154318
$ws.Range("I132").Value = 29413860
$ws.Range("J132").Value = 3556.375
$ws.Range("K132").Value = 88241580
$ws.Range("L132").Value = 10669.125
$ws.Range("M132").Value = -88239050
$ws.Range("N132").Value = -15729.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3775
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 6100
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 6100
$ws.Range("M3").Value = -2886
$ws.Range("N3").Value = -6328

$ws.Range("H75").Value = 11185.6
$ws.Range("I75").Value = 2371.2
$ws.Range("K75").Value = 2371.2
$ws.Range("M75").Value = -1435.2

$ws.Range("H78").Value = 11185.6
$ws.Range("I78").Value = 2371.2
$ws.Range("K78").Value = 7113.599999999999
$ws.Range("M78").Value = -2433.599999999999

$ws.Range("H94").Value = 1878.3
$ws.Range("I94").Value = 1872.875
$ws.Range("K94").Value = 1872.875
$ws.Range("M94").Value = -1421.875

$ws.Range("H105").Value = 3782.4119
$ws.Range("I105").Value = 2282.8572
$ws.Range("J105").Value = 4832.1
$ws.Range("K105").Value = 2282.8572
$ws.Range("L105").Value = 4832.1
$ws.Range("M105").Value = -535.8571999999999
$ws.Range("N105").Value = -8326.1

$ws.Range("H107").Value = 5004.8
$ws.Range("J107").Value = 5253.25
$ws.Range("L107").Value = 5253.25
$ws.Range("N107").Value = -9093.25

$ws.Range("H137").Value = 64592.332
$ws.Range("J137").Value = 64592.332
$ws.Range("L137").Value = 64592.332
$ws.Range("N137").Value = -74792.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2032
$ws.Range("I16").Value = 1827.75
$ws.Range("J16").Value = 2304.3333
$ws.Range("K16").Value = 1827.75
$ws.Range("L16").Value = 2304.3333
$ws.Range("M16").Value = -1540.75
$ws.Range("N16").Value = -2878.3333

$ws.Range("H70").Value = 32902
$ws.Range("J70").Value = 32902
$ws.Range("L70").Value = 32902
$ws.Range("N70").Value = -33532

$ws.Range("H73").Value = 32902
$ws.Range("J73").Value = 32902
$ws.Range("L73").Value = 32902
$ws.Range("N73").Value = -35086

$ws.Range("H105").Value = 3439.0908
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 489.3125
$ws.Range("I107").Value = 364.53845
$ws.Range("J107").Value = 1030
$ws.Range("K107").Value = 364.53845
$ws.Range("L107").Value = 1030
$ws.Range("M107").Value = 1555.46155
$ws.Range("N107").Value = -4870

$ws.Range("H113").Value = 2032
$ws.Range("I113").Value = 1827.75
$ws.Range("J113").Value = 2304.3333
$ws.Range("K113").Value = 1827.75
$ws.Range("L113").Value = 2304.3333
$ws.Range("M113").Value = 342.25
$ws.Range("N113").Value = -6644.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 25303.166
$ws.Range("I18").Value = 27576.182
$ws.Range("K18").Value = 82728.546
$ws.Range("M18").Value = -82559.546

$ws.Range("H120").Value = 1501500
$ws.Range("I120").Value = 1501500
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 4504500
$ws.Range("L120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -4499662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 37900
$ws.Range("J52").Value = 37900
$ws.Range("L52").Value = 37900
$ws.Range("N52").Value = -38418

$ws.Range("H97").Value = 4036.2727
$ws.Range("I97").Value = 2513.2144
$ws.Range("K97").Value = 2513.2144
$ws.Range("M97").Value = -2017.2144

$ws.Range("H113").Value = 1746.0769
$ws.Range("I113").Value = 1833.3334
$ws.Range("J113").Value = 1549.75
$ws.Range("K113").Value = 1833.3334
$ws.Range("L113").Value = 1549.75
$ws.Range("M113").Value = 336.6666
$ws.Range("N113").Value = -5889.75

$ws.Range("H126").Value = 153006
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1339.8
$ws.Range("I22").Value = 1299.5
$ws.Range("J22").Value = 1366.6666
$ws.Range("K22").Value = 1299.5
$ws.Range("L22").Value = 1366.6666
$ws.Range("M22").Value = -1004.5
$ws.Range("N22").Value = -1956.6666

$ws.Range("H27").Value = 1339.8
$ws.Range("I27").Value = 1299.5
$ws.Range("J27").Value = 1366.6666
$ws.Range("K27").Value = 1299.5
$ws.Range("L27").Value = 1366.6666
$ws.Range("M27").Value = -1192.5
$ws.Range("N27").Value = -1580.6666

$ws.Range("H93").Value = 1126.6
$ws.Range("J93").Value = 1233.6666
$ws.Range("L93").Value = 1233.6666
$ws.Range("N93").Value = -3729.6666

$ws.Range("H132").Value = 3603
$ws.Range("I132").Value = 2660.5557
$ws.Range("J132").Value = 6026.4287
$ws.Range("K132").Value = 7981.6671
$ws.Range("L132").Value = 18079.2861
$ws.Range("M132").Value = -5451.6671
$ws.Range("N132").Value = -23139.2861

$ws.Range("H136").Value = 1494.15
$ws.Range("I136").Value = 1061.2354
$ws.Range("J136").Value = 3947.3333
$ws.Range("K136").Value = 3183.7062
$ws.Range("L136").Value = 11841.9999
$ws.Range("M136").Value = -633.7062000000001
$ws.Range("N136").Value = -16941.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H107").Value = 8333974
$ws.Range("I107").Value = 618.9
$ws.Range("J107").Value = 50000750
$ws.Range("K107").Value = 1856.7
$ws.Range("L107").Value = 150002250
$ws.Range("M107").Value = 63.30000000000018
$ws.Range("N107").Value = -150006090

$ws.Range("H113").Value = 714.8570999999999
$ws.Range("I113").Value = 663.5
$ws.Range("J113").Value = 783.3333
$ws.Range("K113").Value = 1990.5
$ws.Range("L113").Value = 2349.9999
$ws.Range("M113").Value = 179.5
$ws.Range("N113").Value = -6689.9999

$ws.Range("H132").Value = 1548.3235
$ws.Range("I132").Value = 1168
$ws.Range("J132").Value = 2245.5833
$ws.Range("K132").Value = 3504
$ws.Range("L132").Value = 6736.749899999999
$ws.Range("M132").Value = -974
$ws.Range("N132").Value = -11796.7499
